# Update average_county_temperature (column AD) values with NOAA data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AD2:AD25").Value = 21.79166666666666
$ws.Range("AD46:AD56").Value = -1.819444444444444
$ws.Range("AD73:AD76").Value = 21.28240740740739
$ws.Range("AD81:AD108").Value = 12.93898809523811
